# Auto-generated script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while always keeping it as literal text,
# even when it looks like a number (e.g. "1.010") or a percentage string.
# We build the text in a scratch cell via a formula ("=""...""") which forces
# Excel to treat it as a string, copy it, and paste-special values-only into
# the destination so no extra number-format/style gets attached to the cell.
function Set-TextValue([string]$cellAddr, [string]$value) {
    $scratch = $ws.Range("ZZ1000")
    $scratch.Formula = '="' + $value + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

Set-TextValue "D2" '27.830.75'
Set-TextValue "E2" '  +0.48%  '
Set-TextValue "D3" '1.858.85'
Set-TextValue "E3" '  +0.05%  '
Set-TextValue "D4" '1.010'
Set-TextValue "E4" '  -2.16%  '
Set-TextValue "D5" '319.51'
Set-TextValue "D6" '1.009'
Set-TextValue "E6" '  -2.04%  '
Set-TextValue "D7" '0.4305'
Set-TextValue "E7" '  -2.22%  '
Set-TextValue "D8" '0.3754'
Set-TextValue "E8" '  -1.43%  '
Set-TextValue "D9" '0.07353'
Set-TextValue "E9" '  -1.15%  '
Set-TextValue "D10" '0.8784'
Set-TextValue "E10" '  -0.86%  '
Set-TextValue "D11" '21.64'
Set-TextValue "E11" '  +0.07%  '
Set-TextValue "D12" '1.856.67'
Set-TextValue "E12" '  -0.25%  '
Set-TextValue "D13" '6.752'
Set-TextValue "E13" '  +0.07%  '
Set-TextValue "D14" '5.451'
Set-TextValue "E14" '  -1.62%  '
Set-TextValue "D15" '0.07133'
Set-TextValue "E15" '  -0.57%  '
Set-TextValue "D16" '89.14'
Set-TextValue "E16" '  +4.33%  '
Set-TextValue "D17" '1.012'
Set-TextValue "E17" '  -2.22%  '
Set-TextValue "D18" '0.000009010'
Set-TextValue "E18" '  -1.03%  '
Set-TextValue "E19" '  -1.85%  '
Set-TextValue "D20" '15.46'
Set-TextValue "E20" '  -0.44%  '
Set-TextValue "D21" '27.810.93'
Set-TextValue "E21" '  +0.28%  '
Set-TextValue "D22" '5.218'
Set-TextValue "E22" '  -1.75%  '
Set-TextValue "D23" '11.09'
Set-TextValue "E23" '  -1.78%  '
Set-TextValue "D24" '2.080.99'
Set-TextValue "E24" '  -0.60%  '
Set-TextValue "D25" '1.985'
Set-TextValue "E25" '  -1.91%  '
Set-TextValue "D26" '155.30'
Set-TextValue "E26" '  -1.73%  '
Set-TextValue "D27" '18.67'
Set-TextValue "E27" '  -0.73%  '
Set-TextValue "D28" '2.177'
Set-TextValue "E28" '  +9.60%  '
Set-TextValue "D29" '5.380'
Set-TextValue "E29" '  -0.05%  '
Set-TextValue "D30" '119.12'
Set-TextValue "E30" '  +0.94%  '
Set-TextValue "D31" '0.08947'
Set-TextValue "E31" '  -0.82%  '
Set-TextValue "D32" '1.233'
Set-TextValue "E32" '  +1.32%  '
Set-TextValue "D33" '0.7803'
Set-TextValue "E33" '  -0.10%  '
Set-TextValue "D34" '4.563'
Set-TextValue "E34" '  -0.33%  '
Set-TextValue "D35" '2.928'
Set-TextValue "E35" '  -2.45%  '
Set-TextValue "D36" '1.010'
Set-TextValue "E36" '  -2.03%  '
Set-TextValue "D37" '1.131'
Set-TextValue "E37" '  -1.43%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D38" '0.05359'
Set-TextValue "E38" '  +1.37%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D39" '0.01985'
Set-TextValue "E39" '  +0.23%  '
Set-TextValue "D40" '7.280'
Set-TextValue "E40" '  +5.82%  '
Set-TextValue "D41" '2.908'
Set-TextValue "E41" '  +1.64%  '
Set-TextValue "D42" '0.1698'
Set-TextValue "E42" '  +0.65%  '
Set-TextValue "D43" '0.5141'
Set-TextValue "E43" '  -1.29%  '
Set-TextValue "D44" '8.847'
Set-TextValue "E44" '  -0.32%  '
Set-TextValue "D45" '10.75'
Set-TextValue "E45" '  +0.71%  '
Set-TextValue "D46" '108.31'
Set-TextValue "E46" '  -2.07%  '
Set-TextValue "D47" '0.4779'
Set-TextValue "E47" '  +1.09%  '
Set-TextValue "D48" '0.06481'
Set-TextValue "E48" '  -1.89%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D49" '1.693'
Set-TextValue "E49" '  -1.47%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D50" '1.010'
Set-TextValue "E50" '  -2.18%  '
Set-TextValue "D51" '1.848'
Set-TextValue "E51" '  -2.94%  '

# Clean up the scratch cell and copy mode
$ws.Range("ZZ1000").Clear()
$excel.CutCopyMode = 0

